# Update the "Results" (column G) values on the "Units" sheet to reflect
# finalized/cleaned config values for the Env Rep team.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Units")

$updates = @{
    6   = 218
    13  = 1063944
    14  = 1063944
    23  = 1056420
    25  = 2328
    26  = 18600
    29  = 18600
    45  = 110738
    47  = 163263
    50  = 163263
    51  = 110738
    66  = 1056420
    67  = 1056420
    76  = 1063944
    99  = 163263
    100 = 1686
    104 = 1063944
    105 = 2328
    113 = 1063944
    114 = 1686
    117 = 1063944
    118 = 1063944
    120 = 1063944
    123 = 18600
    127 = 218
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
